$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply weekly price/volume/date corrections (row-wise shuffle of Fecha/Volumen/Precio fields)
# Row 2
$ws.Range("D2").Value = 44315
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("S2").Value = 1025

# Row 3
$ws.Range("D3").Value = 44410
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("S3").Value = 1025

# Row 4
$ws.Range("D4").Value = 44462
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 19500
$ws.Range("P4").Value = 19750
$ws.Range("S4").Value = 988

# Row 5
$ws.Range("D5").Value = 44442
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("S5").Value = 1025

# Row 6
$ws.Range("D6").Value = 44335
$ws.Range("N6").Value = 19000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 19500
$ws.Range("S6").Value = 975

# Row 7
$ws.Range("D7").Value = 44473
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 19500
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19750
$ws.Range("S7").Value = 988

# Row 8
$ws.Range("D8").Value = 44418
$ws.Range("M8").Value = 200

# Row 9
$ws.Range("D9").Value = 44326
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 19500
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 19750
$ws.Range("S9").Value = 988

# Row 10
$ws.Range("D10").Value = 44474
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("S10").Value = 975

# Row 11
$ws.Range("D11").Value = 44466
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("S11").Value = 1025

# Row 12
$ws.Range("D12").Value = 44448
$ws.Range("M12").Value = 100

# Row 14
$ws.Range("D14").Value = 44407
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 20500
$ws.Range("S14").Value = 1025

# Row 15
$ws.Range("D15").Value = 44350
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("S15").Value = 975

# Row 16
$ws.Range("D16").Value = 44445
$ws.Range("M16").Value = 160
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20500
$ws.Range("S16").Value = 1025

# Row 17
$ws.Range("D17").Value = 44435
$ws.Range("M17").Value = 260
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21115
$ws.Range("S17").Value = 1056

# Row 18
$ws.Range("D18").Value = 44343
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 19500
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19750
$ws.Range("S18").Value = 988

# Row 19
$ws.Range("D19").Value = 44364
$ws.Range("M19").Value = 140

# Row 20
$ws.Range("D20").Value = 44428
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("S20").Value = 1025

# Row 21
$ws.Range("D21").Value = 44333
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 19500
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19750
$ws.Range("S21").Value = 988

# Row 22
$ws.Range("D22").Value = 44431
$ws.Range("M22").Value = 160
$ws.Range("N22").Value = 21000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 21500
$ws.Range("S22").Value = 1075

# Row 23
$ws.Range("D23").Value = 44420
$ws.Range("M23").Value = 160

# Row 24
$ws.Range("D24").Value = 44365
$ws.Range("M24").Value = 100

# Row 25
$ws.Range("D25").Value = 44417
$ws.Range("M25").Value = 160

# Row 26
$ws.Range("D26").Value = 44427
$ws.Range("M26").Value = 200

# Row 27
$ws.Range("D27").Value = 44441
$ws.Range("M27").Value = 160

# Row 28
$ws.Range("D28").Value = 44434
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 20000
$ws.Range("O28").Value = 21000
$ws.Range("P28").Value = 20500
$ws.Range("S28").Value = 1025

# Row 29
$ws.Range("D29").Value = 44301
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = 18000
$ws.Range("O29").Value = 19000
$ws.Range("P29").Value = 18500
$ws.Range("S29").Value = 925

# Row 30
$ws.Range("D30").Value = 44336
